# Actualización automática 2025-09-15 08:55:09
#
# Applies the monthly sales-figures refresh to the three report sheets:
#   - "VENTAS POR GRUPO"     (per-group sales totals)
#   - "VENTA MENSUAL"        (per-month sales totals)
#   - "CUMPLIMIENTO MENSUAL" (budget-vs-sales compliance summary)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M29").Value = 2090.65
$wsGrupo.Range("L31").Value = 760.3200000000001
$wsGrupo.Range("D44").Value = 457.92

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F29").Value = 8733.540000000001
$wsMensual.Range("F31").Value = 3958.26
$wsMensual.Range("F44").Value = 1167.13
$wsMensual.Range("F57").Value = 35029.42

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 240X80 PORCELANATO (row 3)
$wsCumpl.Range("D3").Value = 2635.59
$wsCumpl.Range("E3").Value = 17751.8874217135
$wsCumpl.Range("F3").Value = 0.1292749439022303

# INODOROS (row 6)
$wsCumpl.Range("D6").Value = 474.66
$wsCumpl.Range("E6").Value = 1851.40694516821
$wsCumpl.Range("F6").Value = 0.2040611947932027

# LAVABOS (row 7)
$wsCumpl.Range("D7").Value = 162.55
$wsCumpl.Range("E7").Value = 724.1610162875741
$wsCumpl.Range("F7").Value = 0.1833178984068047

# PIEDRA SINTERIZADA (row 11)
$wsCumpl.Range("D11").Value = 4165.57
$wsCumpl.Range("E11").Value = 15407.4902492497
$wsCumpl.Range("F11").Value = 0.2128216000438501

# TOTAL (row 15)
$wsCumpl.Range("D15").Value = 40676.87
$wsCumpl.Range("E15").Value = 58221.12992509275
$wsCumpl.Range("F15").Value = 0.4113012399725924

# Column F narrowed slightly (25 -> 24 characters) as part of the refresh.
$wsCumpl.Columns.Item(6).ColumnWidth = 23.1666666666667
